function Set-TextValue($cell, $text) {
    $cell.Value = $text
    $v = $cell.Value2
    if ($v -is [double]) {
        $cell.Value = "'" + $text
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range("D2") "56.897.42"
Set-TextValue $ws.Range("E2") "  -0.60%  "

Set-TextValue $ws.Range("D3") "2.974.12"
Set-TextValue $ws.Range("E3") "  -1.48%  "

Set-TextValue $ws.Range("D4") "0.998"
Set-TextValue $ws.Range("E4") "  -0.15%  "

Set-TextValue $ws.Range("D5") "499.90"
Set-TextValue $ws.Range("E5") "  -3.76%  "

Set-TextValue $ws.Range("D6") "137.58"
Set-TextValue $ws.Range("E6") "  -2.58%  "

Set-TextValue $ws.Range("D7") "0.999"
Set-TextValue $ws.Range("E7") "  -0.17%  "

Set-TextValue $ws.Range("D8") "0.429"
Set-TextValue $ws.Range("E8") "  -2.11%  "

Set-TextValue $ws.Range("D9") "7.32"
Set-TextValue $ws.Range("E9") "  -3.53%  "

Set-TextValue $ws.Range("E10") "  -1.35%  "

Set-TextValue $ws.Range("D11") "0.358"
Set-TextValue $ws.Range("E11") "  -0.34%  "

Set-TextValue $ws.Range("D12") "3.472.07"
Set-TextValue $ws.Range("E12") "  -1.90%  "

Set-TextValue $ws.Range("E13") "  -1.95%  "

Set-TextValue $ws.Range("D14") "26.01"
Set-TextValue $ws.Range("E14") "  +0.23%  "

Set-TextValue $ws.Range("D15") "0.0000159"
Set-TextValue $ws.Range("E15") "  -0.23%  "

Set-TextValue $ws.Range("D16") "56.952.55"
Set-TextValue $ws.Range("E16") "  -0.54%  "

Set-TextValue $ws.Range("D17") "6.04"
Set-TextValue $ws.Range("E17") "  +0.41%  "

Set-TextValue $ws.Range("D18") "2.978.87"
Set-TextValue $ws.Range("E18") "  -1.35%  "

Set-TextValue $ws.Range("D19") "12.60"
Set-TextValue $ws.Range("E19") "  -0.74%  "

Set-TextValue $ws.Range("D20") "7.87"
Set-TextValue $ws.Range("E20") "  -1.19%  "

Set-TextValue $ws.Range("D21") "320.95"
Set-TextValue $ws.Range("E21") "  -3.19%  "

Set-TextValue $ws.Range("D22") "1.00"
Set-TextValue $ws.Range("E22") "  +0.10%  "

Set-TextValue $ws.Range("D23") "5.70"
Set-TextValue $ws.Range("E23") "  +0.11%  "

Set-TextValue $ws.Range("D24") "0.488"
Set-TextValue $ws.Range("E24") "  -0.19%  "

Set-TextValue $ws.Range("D25") "63.54"
Set-TextValue $ws.Range("E25") "  -1.24%  "

Set-TextValue $ws.Range("E26") "  +0.48%  "

Set-TextValue $ws.Range("D27") "0.164"
Set-TextValue $ws.Range("E27") "  -4.80%  "

Set-TextValue $ws.Range("D28") "0.0₃0892"
Set-TextValue $ws.Range("E28") "  -3.45%  "

Set-TextValue $ws.Range("D29") "6.52"
Set-TextValue $ws.Range("E29") "  -3.60%  "

Set-TextValue $ws.Range("D30") "7.06"
Set-TextValue $ws.Range("E30") "  -1.52%  "

Set-TextValue $ws.Range("E31") "  -3.47%  "

Set-TextValue $ws.Range("D32") "1.16"
Set-TextValue $ws.Range("E32") "  -4.77%  "

Set-TextValue $ws.Range("D33") "20.20"
Set-TextValue $ws.Range("E33") "  -2.74%  "

Set-TextValue $ws.Range("D34") "155.16"
Set-TextValue $ws.Range("E34") "  -2.04%  "

Set-TextValue $ws.Range("D35") "4.60"
Set-TextValue $ws.Range("E35") "  -1.01%  "

Set-TextValue $ws.Range("D36") "5.78"
Set-TextValue $ws.Range("E36") "  +0.13%  "

Set-TextValue $ws.Range("D37") "1.25"
Set-TextValue $ws.Range("E37") "  -4.07%  "

Set-TextValue $ws.Range("D38") "24.15"
Set-TextValue $ws.Range("E38") "  -0.41%  "

Set-TextValue $ws.Range("D39") "0.0667"
Set-TextValue $ws.Range("E39") "  -1.83%  "

Set-TextValue $ws.Range("B40") "OKB"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D40") "37.64"
Set-TextValue $ws.Range("E40") "  +0.56%  "

Set-TextValue $ws.Range("B41") "RenzoRestakedETH"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextValue $ws.Range("D41") "3.003.08"
Set-TextValue $ws.Range("E41") "  -1.61%  "

Set-TextValue $ws.Range("D42") "0.999"
Set-TextValue $ws.Range("E42") "  -0.07%  "

Set-TextValue $ws.Range("D43") "3.73"
Set-TextValue $ws.Range("E43") "  -0.04%  "

Set-TextValue $ws.Range("D44") "0.638"
Set-TextValue $ws.Range("E44") "  -2.34%  "

Set-TextValue $ws.Range("B45") "Stacks"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D45") "1.39"
Set-TextValue $ws.Range("E45") "  -3.41%  "

Set-TextValue $ws.Range("B46") "Maker"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D46") "2.191.85"
Set-TextValue $ws.Range("E46") "  -4.90%  "

Set-TextValue $ws.Range("B47") "ONDO"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws.Range("D47") "0.941"
Set-TextValue $ws.Range("E47") "  -6.62%  "

Set-TextValue $ws.Range("B48") "Cosmos"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D48") "5.94"
Set-TextValue $ws.Range("E48") "  +0.56%  "

Set-TextValue $ws.Range("D49") "0.0235"
Set-TextValue $ws.Range("E49") "  -3.66%  "

Set-TextValue $ws.Range("D50") "19.22"
Set-TextValue $ws.Range("E50") "  -1.51%  "

Set-TextValue $ws.Range("D51") "1.80"
Set-TextValue $ws.Range("E51") "  -10.86%  "
